# Apply odds updates to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 3.7
    "I2"  = 1.9
    "W2"  = 1.29
    "AB2" = 23
    "AC2" = 13

    "G3"  = 2
    "H3"  = 3.2
    "I3"  = 4
    "L3"  = 4.75
    "R3"  = 1.5
    "W3"  = 1.57
    "X3"  = 2.25
    "AE3" = 21
    "AG3" = 6.5
    "AL3" = 8.5
    "AR3" = 1.93
    "AS3" = 1.93

    "G4"  = 1.73
    "I4"  = 5
    "M4"  = 1.07
    "N4"  = 9
    "O4"  = 1.36
    "P4"  = 3.2
    "Q4"  = 2.2
    "R4"  = 1.67
    "U4"  = 4
    "V4"  = 1.25

    "H5"  = 3.8
    "I5"  = 1.62
    "L5"  = 2.2
    "M5"  = 1.03
    "N5"  = 15
    "O5"  = 1.2
    "P5"  = 4.33
    "Q5"  = 1.67
    "R5"  = 2.15
    "U5"  = 2.63
    "V5"  = 1.44
    "W5"  = 1.3
    "AE5" = 41

    "G7"  = 1.48
    "M7"  = 1.03
    "N7"  = 15
    "O7"  = 1.2
    "P7"  = 4.33
    "Q7"  = 1.67
    "R7"  = 2.15
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
